# Apply attendance updates to Sheet1.
# Columns: A=Date, B=Roll, C=Name, D=Total Attendance Count, E=Real,
#          F=Duplicate, G=Invalid, H=Absent
# Rows 3..18 each correspond to one attendance date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# row -> list of columns that flip from 0 to 1
$updates = @{
    3  = @("G", "H")
    4  = @("D", "E")
    5  = @("H")
    6  = @("D", "E")
    7  = @("H")
    8  = @("H")
    9  = @("D", "E")
    10 = @("H")
    11 = @("H")
    12 = @("D", "E")
    13 = @("D", "E")
    14 = @("G", "H")
    15 = @("H")
    16 = @("H")
    17 = @("H")
    18 = @("H")
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row]) {
        $ws.Range("$col$row").Value = 1
    }
}
